# Add new analysis of resazurin (plates 1-19) measured on 2025-09-04
# These mirror rows 2-20 (date 2025-09-03, plate1..plate19) with the same
# plate labels but a new date and (identical) measurement values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 20250904
$startRow = 40

for ($i = 0; $i -lt 19; $i++) {
    $row = $startRow + $i
    $srcRow = 2 + $i

    $ws.Cells.Item($row, 1).Value = $newDate
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($srcRow, 2).Value()
    $ws.Cells.Item($row, 3).Value = $ws.Cells.Item($srcRow, 3).Value()
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($srcRow, 4).Value()
}

# Update the view state to match the author's final scroll/selection position
# (scrolled so row 35 is the top visible row, with C60 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("C60").Select()
